$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values per repull / mean calculation fix
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = -8
$ws.Range("F10").Value = 1
